# Adds a "time_taken" column (F) to the worksheet, mirroring the header
# style of the existing "panel" column (E) and filling each data row
# (2-76) with its corresponding timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 -------------------------------------------------
# Copy the formatting of E1 (bold, centered, bordered header style)
# onto F1, then set its text. Using Copy/PasteSpecial(Formats) reuses
# the existing header cell style instead of creating a new one.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F1").Value = "time_taken"

# --- Data rows 2-76 ("time_taken" timestamps) ------------------------
$timeTaken = @(
    "2021-10-05 13:38:57.191927", "2021-10-05 13:38:57.191940", "2021-10-05 13:38:57.191944",
    "2021-10-05 13:38:57.191947", "2021-10-05 13:38:57.191951", "2021-10-05 13:38:57.191954",
    "2021-10-05 13:38:57.191957", "2021-10-05 13:38:57.191960", "2021-10-05 13:38:57.191963",
    "2021-10-05 13:38:57.191966", "2021-10-05 13:38:57.191969", "2021-10-05 13:38:57.191972",
    "2021-10-05 13:38:57.191975", "2021-10-05 13:38:57.191978", "2021-10-05 13:38:57.191981",
    "2021-10-05 13:38:57.191984", "2021-10-05 13:38:57.191988", "2021-10-05 13:38:57.191991",
    "2021-10-05 13:38:57.191994", "2021-10-05 13:38:57.191997", "2021-10-05 13:38:57.191999",
    "2021-10-05 13:38:57.192002", "2021-10-05 13:38:57.192005", "2021-10-05 13:38:57.192008",
    "2021-10-05 13:38:57.192012", "2021-10-05 13:38:57.192015", "2021-10-05 13:38:57.192018",
    "2021-10-05 13:38:57.192021", "2021-10-05 13:38:57.192023", "2021-10-05 13:38:57.192026",
    "2021-10-05 13:38:57.192029", "2021-10-05 13:38:57.192032", "2021-10-05 13:38:57.192035",
    "2021-10-05 13:38:57.192039", "2021-10-05 13:38:57.192041", "2021-10-05 13:38:57.192044",
    "2021-10-05 13:38:57.192047", "2021-10-05 13:38:57.192050", "2021-10-05 13:38:57.192053",
    "2021-10-05 13:38:57.192056", "2021-10-05 13:38:57.192059", "2021-10-05 13:38:57.192062",
    "2021-10-05 13:38:57.192065", "2021-10-05 13:38:57.192068", "2021-10-05 13:38:57.192071",
    "2021-10-05 13:38:57.192074", "2021-10-05 13:38:57.192077", "2021-10-05 13:38:57.192080",
    "2021-10-05 13:38:57.192083", "2021-10-05 13:38:57.192086", "2021-10-05 13:38:57.192089",
    "2021-10-05 13:38:57.192092", "2021-10-05 13:38:57.192096", "2021-10-05 13:38:57.192099",
    "2021-10-05 13:38:57.192102", "2021-10-05 13:38:57.192105", "2021-10-05 13:38:57.192108",
    "2021-10-05 13:38:57.192111", "2021-10-05 13:38:57.192113", "2021-10-05 13:38:57.192116",
    "2021-10-05 13:38:57.192119", "2021-10-05 13:38:57.192122", "2021-10-05 13:38:57.192125",
    "2021-10-05 13:38:57.192128", "2021-10-05 13:38:57.192133", "2021-10-05 13:38:57.192136",
    "2021-10-05 13:38:57.192139", "2021-10-05 13:38:57.192142", "2021-10-05 13:38:57.192145",
    "2021-10-05 13:38:57.192148", "2021-10-05 13:38:57.192151", "2021-10-05 13:38:57.192154",
    "2021-10-05 13:38:57.192157", "2021-10-05 13:38:57.192160", "2021-10-05 13:38:57.192163"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
